$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (columns B:E): 1,2,3,4 -> 16,20,16,20
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (Subj "CON") updated measurement values for columns B:E
$ws.Range("B2").Value = 387.3451681125806
$ws.Range("C2").Value = 458.35405600325987
$ws.Range("D2").Value = 383.1533072526862
$ws.Range("E2").Value = 461.3256402203696

# Row 3 (Subj "STR") updated measurement values for columns B:E
$ws.Range("B3").Value = 387.1734780318735
$ws.Range("C3").Value = 473.65284788893786
$ws.Range("D3").Value = 392.6553446345344
$ws.Range("E3").Value = 464.74641813881294

# Update the saved selection to match the edited range
$ws.Range("B1:E3").Select() | Out-Null
